$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-07-01"

# Update the header label (shared string) for the "through" column
$ws.Range("I1").Value = "2022 (through 07-01)"

# Add the new July data point in the 2022 column
$ws.Range("I8").Value = 1

# Update the Total row for the 2022 column
$ws.Range("I14").Value = 807
